$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark (it used to sit at the very end of
#    the last paragraph). It will be re-created further down, right after
#    the rewritten "Description" paragraph.
# ----------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ----------------------------------------------------------------------
# 2) Rewrite the "Description" paragraph (paragraph 4) with the new,
#    longer text describing both the classifier and the web interface.
# ----------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

$newDescription = "I will train" + `
    " a document classifier to detect " + $openQuote + "types" + $closeQuote + " of partisan news using long" + `
    " term - " + `
    "short term memory neural networks. " + `
    "Then" + `
    ", " + `
    "I will " + `
    "create a web interface" + `
    " for human usability" + `
    " that uses " + `
    "a form to send article text to " + `
    "the trained neural network model" + `
    ", and returns a graphical representation of the results. "

$r4.Text = $newDescription

# ----------------------------------------------------------------------
# 3) Delete the blank paragraph that used to follow the description.
# ----------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5.Range.Delete()

# ----------------------------------------------------------------------
# 4) Re-add the "_GoBack" bookmark right at the end of the description
#    paragraph (collapsed range). A temporary marker character is used
#    to work around collapsed-range placement right before a paragraph
#    mark, then removed again once the bookmark is anchored.
# ----------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range
$r4.End = $r4.End - 1
$r4.InsertAfter("~")
$bmPos = $r4.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
$markerRange = $d.Range($bmPos, $bmPos + 1)
$markerRange.Delete()

# ----------------------------------------------------------------------
# 5) Trim the "Scrape my dataset" milestone: drop the "For robustness,
#    also try BBC, PBS, Baribart, fiveThirtyEight, etc." aside.
# ----------------------------------------------------------------------
$scrapePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*Scrape my dataset*") {
        $scrapePara = $candidate
        break
    }
}

$searchRange = $scrapePara.Range.Duplicate
$found = $searchRange.Find.Execute("For robustness, also try BBC", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $paraEnd = $scrapePara.Range.End - 1
    $trimRange = $d.Range($searchRange.Start, $paraEnd)
    $trimRange.Delete()
}
